$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Profiles")

# --- Row 1: clear stray empty cells in B1/C1, keep A1 title ---
$ws.Range("B1").Value = $null
$ws.Range("C1").Value = $null

# --- Row 2: replace "Date export" / date columns with a single generated-on note ---
$ws.Range("A2").Value = "Généré le 2026-01-16 09:00:00"
$ws.Range("B2").Value = $null
$ws.Range("C2").Value = $null

# --- Row 4: normalize/adjust country + channel values ---
$ws.Range("B4").Value = " TOGO "
$ws.Range("C4").Value = "Mobile"

# --- Row 5: country casing tweak (channel WEB untouched) ---
$ws.Range("B5").Value = "Benin"

# --- Row 6: country spacing + channel filled in ---
$ws.Range("B6").Value = "  Togo"
$ws.Range("C6").Value = " web"

# --- Row 7: country trailing space, channel cleared ---
$ws.Range("B7").Value = "Togo "
$ws.Range("C7").Value = ""

# --- Row 8: country trimmed, channel recased ---
$ws.Range("B8").Value = "Bénin"
$ws.Range("C8").Value = "Web"

# --- Row 9: country upcased, channel filled in ---
$ws.Range("B9").Value = "TOGO"
$ws.Range("C9").Value = "mobile"

# --- Row 10: channel re-labelled ---
$ws.Range("C10").Value = "Mobile "

# --- Row 11: brand new record appended ---
$ws.Range("A11").Value = "u008"
$ws.Range("B11").Value = " benin"
$ws.Range("C11").Value = " web"

# --- Sheet2 ("Notes"): collapse the two free-text lines into one shorter note ---
$ws2 = $wb.Worksheets.Item("Notes")
$ws2.Range("A2").Value = "Feuille à ignorer (texte libre)."
$ws2.Range("A3").Value = $null
